$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend the table with a new "2023" column (K), mirroring column J's
# formatting (same number format / borders / alignment as the rest of the
# row), then fill in the new year's figures.
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 871.3
$ws.Range("K5").Value = 485.6
$ws.Range("K6").Value = 1010.7
